$d = $word.ActiveDocument

# Remember where the original (old) content ends before we insert anything.
$oldEnd = $d.Paragraphs(2).Range.End

# Insert the replacement paragraph (merged text, line break, simplified
# formatting) immediately before the existing content.
$insertionPoint = $d.Range(0, 0)
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>Im Namen Gottes, des Vaters und des Sohnes und des Heiligen Geistes.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/></w:rPr><w:t>Gemeinde: Amen.</w:t></w:r></w:p>'
$insertionPoint.InsertXML($newParaXml)

# The old two paragraphs now sit right after the newly inserted one; the
# document grew by the length of the inserted XML's text, so recompute the
# deletion range from the current document length.
$newEnd = $d.Content.End
$oldRange = $d.Range($newEnd - $oldEnd, $newEnd)
$oldRange.Delete()
